# gov pov results working
#
# Updates the distribution-percentage formulas on the MHDV_distribution and
# LDV_distribution sheets to reflect new totals/counts, and refreshes the
# active-sheet/selection state (MHDV -> LDV_distribution becomes the active
# tab, with new selections on several sheets).

$wb = $excel.ActiveWorkbook

# --- MHDV: selection moves from K17 to D2:D21, tab no longer selected ---
$ws1 = $wb.Worksheets.Item("MHDV")
$ws1.Activate()
$ws1.Range("D2:D21").Select()

# --- MHDV_distribution: updated denominators/numerators + new selection ---
$ws2 = $wb.Worksheets.Item("MHDV_distribution")
$ws2.Range("C2").Formula = "=(24/(511-B2))*100"
$ws2.Range("D2").Formula = "=(183/(511-B2))*100"
$ws2.Range("E2").Formula = "=(302/(511-B2))*100"
$ws2.Range("C3").Formula = "=(32/(654-B3))*100"
$ws2.Range("D3").Formula = "=(149/(654-B3))*100"
$ws2.Range("E3").Formula = "=(475/(654-B3))*100"
$ws2.Range("C4").Formula = "=(0/(41-B4))*100"
$ws2.Range("D4").Formula = "=(18/(41-B4))*100"
$ws2.Range("E4").Formula = "=(21/(41-B4))*100"
$ws2.Activate()
$ws2.Range("D5").Select()

# --- LDV_distribution: updated denominators/numerators, becomes active tab ---
$ws4 = $wb.Worksheets.Item("LDV_distribution")
$ws4.Range("C2").Formula = "=(307/(616-B2))*100"
$ws4.Range("D2").Formula = "=(127/(616-B2))*100"
$ws4.Range("E2").Formula = "=(182/(616-B2))*100"
$ws4.Range("C4").Formula = "=(21/(45-B4))*100"
$ws4.Range("D4").Formula = "=(5/(45-B4))*100"
$ws4.Range("E4").Formula = "=(19/(45-B4))*100"
$ws4.Activate()
$ws4.Range("F10").Select()
